$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header labels to align with new Excel file structure
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "IDAM Roles"

# Reset the active cell selection to C1
$ws.Range("C1").Select()
